$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 28
$ws.Range("C2").Value = "face/face065.png"
$ws.Range("D2").Value = "hoffen"
$ws.Range("B3").Value = 42
$ws.Range("C3").Value = "car/car082.png"
$ws.Range("D3").Value = "jubeln"
$ws.Range("B4").Value = 126
$ws.Range("C4").Value = "car/car116.png"
$ws.Range("D4").Value = "bleiben"
$ws.Range("B5").Value = 104
$ws.Range("C5").Value = "car/car094.png"
$ws.Range("D5").Value = "lehnen"
$ws.Range("B6").Value = 36
$ws.Range("C6").Value = "car/car064.png"
$ws.Range("D6").Value = "starten"
$ws.Range("B7").Value = 81
$ws.Range("C7").Value = "face/face106.png"
$ws.Range("D7").Value = "hauen"
$ws.Range("B8").Value = 15
$ws.Range("C8").Value = "face/face085.png"
$ws.Range("D8").Value = "gründen"
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = "car/car090.png"
$ws.Range("D9").Value = "bitten"
$ws.Range("B10").Value = 21
$ws.Range("C10").Value = "car/car072.png"
$ws.Range("D10").Value = "schenken"
$ws.Range("B11").Value = 112
$ws.Range("C11").Value = "face/face123.png"
$ws.Range("D11").Value = "füllen"
$ws.Range("B12").Value = 120
$ws.Range("C12").Value = "car/car123.png"
$ws.Range("D12").Value = "hupen"
$ws.Range("B13").Value = 87
$ws.Range("C13").Value = "face/face120.png"
$ws.Range("D13").Value = "drohen"
$ws.Range("B14").Value = 35
$ws.Range("C14").Value = "car/car100.png"
$ws.Range("D14").Value = "antun"
$ws.Range("B15").Value = 109
$ws.Range("C15").Value = "car/car081.png"
$ws.Range("D15").Value = "backen"
$ws.Range("B16").Value = 105
$ws.Range("C16").Value = "car/car093.png"
$ws.Range("D16").Value = "strahlen"
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = "car/car073.png"
$ws.Range("D17").Value = "klappen"
$ws.Range("B18").Value = 113
$ws.Range("C18").Value = "car/car104.png"
$ws.Range("D18").Value = "schätzen"
$ws.Range("B19").Value = 8
$ws.Range("C19").Value = "face/face094.png"
$ws.Range("D19").Value = "rücken"
$ws.Range("B20").Value = 76
$ws.Range("C20").Value = "face/face068.png"
$ws.Range("D20").Value = "schmecken"
$ws.Range("B21").Value = 45
$ws.Range("C21").Value = "car/car096.png"
$ws.Range("D21").Value = "krachen"
$ws.Range("B22").Value = 63
$ws.Range("C22").Value = "car/car074.png"
$ws.Range("D22").Value = "husten"
$ws.Range("B23").Value = 34
$ws.Range("C23").Value = "face/face075.png"
$ws.Range("D23").Value = "töten"
$ws.Range("B24").Value = 61
$ws.Range("C24").Value = "face/face101.png"
$ws.Range("D24").Value = "sieben"
$ws.Range("B25").Value = 60
$ws.Range("C25").Value = "car/car091.png"
$ws.Range("D25").Value = "raten"
$ws.Range("B26").Value = 56
$ws.Range("C26").Value = "face/face099.png"
$ws.Range("D26").Value = "saufen"
$ws.Range("B27").Value = 57
$ws.Range("C27").Value = "face/face107.png"
$ws.Range("D27").Value = "wenden"
$ws.Range("B28").Value = 47
$ws.Range("C28").Value = "face/face091.png"
$ws.Range("D28").Value = "liefern"
$ws.Range("B29").Value = 86
$ws.Range("C29").Value = "face/face067.png"
$ws.Range("D29").Value = "dauern"
$ws.Range("B30").Value = 41
$ws.Range("C30").Value = "face/face078.png"
$ws.Range("D30").Value = "regnen"
$ws.Range("B31").Value = 102
$ws.Range("C31").Value = "face/face084.png"
$ws.Range("D31").Value = "langen"
$ws.Range("B32").Value = 13
$ws.Range("C32").Value = "face/face095.png"
$ws.Range("D32").Value = "mieten"
$ws.Range("B33").Value = 74
$ws.Range("C33").Value = "car/car097.png"
$ws.Range("D33").Value = "stärken"
